$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Turn F5:F11 into one shared formula (SUM(D*E) per row) ---
$ws.Range("F5:F11").Formula = "=SUM(D5*E5)"

# --- Row 16: add total formula (was blank) ---
$ws.Range("F16").Formula = "=SUM(D16*E16)"

# --- Row 19: add total formula and pick up row16's style (s=6) for F19 ---
$ws.Range("F16").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("F19").Formula = "=SUM(D19*E19)"

# --- Insert a new row 20 (shifts old row 21 "Total costs" down to row 22) ---
$ws.Rows("20").Insert() | Out-Null

# Give the new row 20 cells A:E the same formatting as row 19, and F20 the
# same formatting as F16, then fill in the values (Lazy Susan / Bylim).
$ws.Range("A19:E19").Copy() | Out-Null
$ws.Range("A20:E20").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("F16").Copy() | Out-Null
$ws.Range("F20").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("A20").Value = "Lazy Susan"
$ws.Range("B20").Value = "Bylim"
$ws.Range("D20").Value = 3
$ws.Range("F20").Value = 65

# --- Fix the grand-total formula (now on row 22) to include row 20 ---
$ws.Range("F22").Formula = "=SUM(F5:F20)"

# --- Row 7: new microcontroller line (Kingduino) ---
$ws.Range("B7").Value = "Kingduino"
$ws.Range("E7").Value = 7.81
$ws.Range("H7").Value = "https://hobbyking.com/nl_nl/kingduino-atmel-atmega-328-pu.html"

# --- Restore the selection the author left the sheet on ---
$ws.Range("L18").Select() | Out-Null
